# releaseInneripDialog.xlsx — update per commit "Upload files to console based on cn181107"
#
# Net content change: the dialog's "OK" confirmation text grows from
# "Confirm to release private IP" to "Confirm to release this private IP".
# Also widen column C to fit the longer caption, move the active selection,
# and set the sheet's page setup (paper size / orientation) as captured by
# the reviewer's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the English confirmation string shown in C2.
$ws.Range("C2").Value = "Confirm to release this private IP"

# 2. Column C needs to be wide enough for the new text.
$ws.Columns.Item(3).ColumnWidth = 28.17

# 3. Page setup picked up by the reviewer's save (A4/Letter-ish single sheet,
#    portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# 4. Leave the cursor where the reviewer left it.
[void]$ws.Range("C17").Select()
